$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: 3000 -> 3200
$ws.Range("B2").Value = 3200

# Row 6: was "Ck2 White Diamond " / 4500 / 179
#        now "Tv10 C4 (Готовый)" / 6000 / 253, with custom row height 13.5
$ws.Range("A6").Value = "Tv10 C4 (Готовый)"
$ws.Range("B6").Value = 6000
$ws.Range("C6").Value = 253
$ws.Rows.Item(6).RowHeight = 13.5

# Row 7: was "Tv10 C4 (Готовый)" / 6000 / 253 (custom row height 13.5)
#        now "Ck05 White Diamond " / 1500 / 235, default row height
$ws.Range("A7").Value = "Ck05 White Diamond "
$ws.Range("B7").Value = 1500
$ws.Range("C7").Value = 235
$ws.Rows.Item(7).AutoFit()

# Update selection to match target view state
$ws.Range("F14").Select()
